# Apply "add programing for linked list problem" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Easy")

# Set F2, F4, F5 to "Javascript", matching the style already used for the
# same column/value combination elsewhere in the sheet (e.g. F14:F16).
$template = $ws.Range("F14")
$template.Copy()
$cells = @("F2", "F4", "F5")
foreach ($addr in $cells) {
    $ws.Range($addr).PasteSpecial(-4122) # xlPasteFormats
    $ws.Range($addr).Value = "Javascript"
}

# Update the active selection to F2 (also resets the scrolled-to cell).
$ws.Range("F2").Select()
